$wb = $excel.ActiveWorkbook

# --- Update "want to go" counts (F column) across all sheets ---

# Sheet: 展览
$ws = $wb.Worksheets.Item('展览')
$ws.Range('F3').Value = 1232
$ws.Range('F4').Value = 1294
$ws.Range('F6').Value = 181
$ws.Range('F7').Value = 560
$ws.Range('F8').Value = 26
$ws.Range('F9').Value = 354
$ws.Range('F11').Value = 1275
$ws.Range('F12').Value = 29362
$ws.Range('F13').Value = 4560
$ws.Range('F14').Value = 47
$ws.Range('F17').Value = 54
$ws.Range('F18').Value = 43
$ws.Range('F19').Value = 22
$ws.Range('F20').Value = 29
$ws.Range('F21').Value = 351
$ws.Range('F22').Value = 10
$ws.Range('F23').Value = 646
$ws.Range('F24').Value = 281
$ws.Range('F25').Value = 292
$ws.Range('F26').Value = 362
$ws.Range('F28').Value = 84
$ws.Range('F29').Value = 8
$ws.Range('F30').Value = 673
$ws.Range('F31').Value = 221
$ws.Range('F32').Value = 106
$ws.Range('F33').Value = 561
$ws.Range('F34').Value = 82
$ws.Range('F36').Value = 652
$ws.Range('F37').Value = 247
$ws.Range('F38').Value = 44
$ws.Range('F39').Value = 9

# Sheet: 演出
$ws = $wb.Worksheets.Item('演出')
$ws.Range('F6').Value = 389
$ws.Range('F7').Value = 916
$ws.Range('F8').Value = 2
$ws.Range('F10').Value = 91
$ws.Range('F11').Value = 277
$ws.Range('F12').Value = 4267
$ws.Range('F17').Value = 49
$ws.Range('F18').Value = 4
$ws.Range('F23').Value = 4255

# Sheet: 本地生活
$ws = $wb.Worksheets.Item('本地生活')
$ws.Range('F2').Value = 301
$ws.Range('F3').Value = 270
$ws.Range('F4').Value = 1229
$ws.Range('F5').Value = 305

# Sheet: 全部类型
$ws = $wb.Worksheets.Item('全部类型')
$ws.Range('F2').Value = 301
$ws.Range('F3').Value = 270
$ws.Range('F4').Value = 1229
$ws.Range('F7').Value = 389
$ws.Range('F8').Value = 305
$ws.Range('F9').Value = 916
$ws.Range('F10').Value = 1232
$ws.Range('F11').Value = 1294
$ws.Range('F12').Value = 181
$ws.Range('F13').Value = 560
$ws.Range('F14').Value = 26
$ws.Range('F15').Value = 354
$ws.Range('F18').Value = 1275
$ws.Range('F20').Value = 91
$ws.Range('F26').Value = 49
$ws.Range('F27').Value = 49
$ws.Range('F29').Value = 54
$ws.Range('F30').Value = 22
$ws.Range('F31').Value = 4
$ws.Range('F32').Value = 29
$ws.Range('F34').Value = 351
$ws.Range('F35').Value = 646
$ws.Range('F36').Value = 281
$ws.Range('F38').Value = 84
$ws.Range('F39').Value = 8
$ws.Range('F40').Value = 673
$ws.Range('F42').Value = 221
$ws.Range('F43').Value = 106
$ws.Range('F46').Value = 82
$ws.Range('F48').Value = 652
$ws.Range('F49').Value = 247
$ws.Range('F50').Value = 44

# --- Sheet "全部类型": row content updates (rows 19, 21, 22) ---
# A new exhibition entry ("萤火虫动漫游戏嘉年华") is inserted at row 19
# (sorted by date before the existing "冰兔2024" entry), which pushes the
# text of the following rows down by one; the original row 22 content
# ("昨日重现") is dropped from the merged view.
$ws = $wb.Worksheets.Item('全部类型')

# Row 19
$ws.Range('B19').NumberFormat = "@"
$ws.Range('B19').Value = '2024-07-19'
$ws.Range('B19').Style = "Normal"
$ws.Range('C19').Value = '广州·萤火虫动漫游戏嘉年华 × KKWORLD2024 快看漫画乐园'
$ws.Range('D19').Value = '新港东路1000号 保利世贸博览馆'
$ws.Range('E19').Value = '2024.07.19 09:00-07.22 17:00'
$ws.Range('F19').Value = 29362
$ws.Range('G19').Value = '已售罄'
$ws.Range('H19').Value = 'https://show.bilibili.com/platform/detail.html?id=87210'
$ws.Range('I19').Value = '//i1.hdslb.com/bfs/openplatform/202406/DTCdOTPs1718177177472.jpeg'

# Row 21
$ws.Range('C21').Value = '广州·冰兔2024线下live「过去和未来」'
$ws.Range('D21').Value = '恩宁路265号三层四层自编01 MAO Livehouse广州（永庆坊店）'
$ws.Range('E21').Value = '2024.07.20 20:00-07.20 22:00'
$ws.Range('F21').Value = 91
$ws.Range('G21').Value = 198
$ws.Range('H21').Value = 'https://show.bilibili.com/platform/detail.html?id=87546'
$ws.Range('I21').Value = '//i2.hdslb.com/bfs/openplatform/202406/2X09PE1a1718611339266.jpeg'

# Row 22
$ws.Range('B22').NumberFormat = "@"
$ws.Range('B22').Value = '2024-07-20'
$ws.Range('B22').Style = "Normal"
$ws.Range('C22').Value = '广州·跨越二次元ACG神级动漫世界巡回演唱会'
$ws.Range('D22').Value = '广州市荔湾区十甫路125号(上下九步行街内)2层 广州平安大戏院'
$ws.Range('E22').Value = '2024.07.20 19:30-07.20 21:10'
$ws.Range('F22').Value = 277
$ws.Range('G22').Value = 280
$ws.Range('H22').Value = 'https://show.bilibili.com/platform/detail.html?id=85353'
$ws.Range('I22').Value = '//i1.hdslb.com/bfs/openplatform/202405/4gACWbPh1715223804704.jpeg'

